$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '66.413.94'
$ws.Range('D3').Value = '3.565.92'
$ws.Range('E3').Value = '  -4.61%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue 'D5' '582.89'
$ws.Range('E5').Value = '  -5.04%  '
Set-TextValue 'D6' '186.58'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('D7').Value = '3.562.32'
Set-TextValue 'D8' '0.613'
$ws.Range('E8').Value = '  -4.32%  '
$ws.Range('E9').Value = '  -0.06%  '
Set-TextValue 'D10' '0.669'
$ws.Range('E10').Value = '  -7.38%  '
$ws.Range('E11').Value = '  -10.25%  '
Set-TextValue 'D12' '53.14'
$ws.Range('E12').Value = '  -7.85%  '
$ws.Range('E13').Value = '  -10.87%  '
Set-TextValue 'D14' '9.82'
$ws.Range('E14').Value = '  -8.19%  '
$ws.Range('D15').Value = '4.135.52'
$ws.Range('E15').Value = '  -4.57%  '
$ws.Range('D16').Value = '3.570.43'
$ws.Range('E16').Value = '  -4.45%  '
Set-TextValue 'D17' '0.125'
$ws.Range('E17').Value = '  -0.85%  '
Set-TextValue 'D18' '18.35'
$ws.Range('E18').Value = '  -5.43%  '
Set-TextValue 'D19' '12.21'
$ws.Range('E19').Value = '  -6.83%  '
$ws.Range('D20').Value = '66.345.19'
$ws.Range('E20').Value = '  -3.90%  '
$ws.Range('E21').Value = '  -7.85%  '
Set-TextValue 'D22' '394.88'
$ws.Range('E22').Value = '  -4.71%  '
Set-TextValue 'D23' '4.34'
$ws.Range('E23').Value = '  -6.38%  '
Set-TextValue 'D24' '85.96'
$ws.Range('E24').Value = '  -4.12%  '
Set-TextValue 'D25' '11.26'
$ws.Range('E25').Value = '  +1.48%  '
Set-TextValue 'D26' '2.90'
Set-TextValue 'D27' '12.47'
$ws.Range('E27').Value = '  -3.98%  '
Set-TextValue 'D28' '6.06'
$ws.Range('E28').Value = '  -0.03%  '
Set-TextValue 'D29' '3.55'
$ws.Range('E29').Value = '  -6.91%  '
Set-TextValue 'D30' '8.95'
$ws.Range('E30').Value = '  -8.02%  '
Set-TextValue 'D31' '31.11'
$ws.Range('E31').Value = '  -6.92%  '
Set-TextValue 'D32' '7.10'
$ws.Range('E32').Value = '  -4.61%  '
Set-TextValue 'D33' '12.18'
$ws.Range('E33').Value = '  -4.48%  '
Set-TextValue 'D34' '621.06'
$ws.Range('E34').Value = '  -1.08%  '
Set-TextValue 'D35' '0.113'
$ws.Range('E35').Value = '  -8.88%  '
Set-TextValue 'D36' '63.48'
$ws.Range('E36').Value = '  -4.14%  '
Set-TextValue 'D37' '41.43'
$ws.Range('E37').Value = '  -8.22%  '
$ws.Range('E38').Value = '  +0.14%  '
Set-TextValue 'D39' '0.397'
$ws.Range('E39').Value = '  -5.64%  '
$ws.Range('D40').Value = '0.0₃0762'
$ws.Range('E40').Value = '  -9.45%  '
$ws.Range('E41').Value = '  -6.43%  '
Set-TextValue 'D42' '0.998'
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = '3.003.53'
$ws.Range('E43').Value = '  +4.96%  '
$ws.Range('E44').Value = '  -8.11%  '
$ws.Range('E45').Value = '  -4.66%  '
Set-TextValue 'D46' '0.0408'
$ws.Range('E46').Value = '  -8.68%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D47' '0.131'
$ws.Range('E47').Value = '  -7.09%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 'D48' '3.10'
$ws.Range('E48').Value = '  -1.68%  '
Set-TextValue 'D49' '8.57'
$ws.Range('E49').Value = '  -6.98%  '
Set-TextValue 'D50' '137.30'
$ws.Range('E50').Value = '  -3.77%  '
Set-TextValue 'D51' '2.74'
$ws.Range('E51').Value = '  -1.55%  '
